$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: phone column was stored as text ("79174445"); normalize it to
# a real number, matching the rest of the A column. ---
$ws.Cells.Item(17, 1).Value = 79174445

# --- Row 18 (new): payment 79174445 (Cash) 2025-08-18T08:51:56 ---

# A18 keeps the numeric-looking phone number as literal text (like the raw
# un-normalized row 17 used to be), so force text via the "@" number format,
# then drop back to the default/Normal style so no formatting sticks to the
# cell itself.
$cA18 = $ws.Cells.Item(18, 1)
$cA18.NumberFormat = "@"
$cA18.Value = "79174445"
$cA18.Style = "Normal"

# B18 (amount) is blank on this row, stored as an empty text cell — use a
# bare apostrophe to force an empty-string text entry, then reset style.
$cB18 = $ws.Cells.Item(18, 2)
$cB18.Value = "'"
$cB18.Style = "Normal"

$ws.Cells.Item(18, 3).Value = "Cash"
$ws.Cells.Item(18, 4).Value = "2025-08-18T08:51:56"
$ws.Cells.Item(18, 5).Value = 20

# F18 (discount_applied) is blank too.
$cF18 = $ws.Cells.Item(18, 6)
$cF18.Value = "'"
$cF18.Style = "Normal"

$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0
$ws.Cells.Item(18, 9).Value = 20
